$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.394.59'
$ws.Range('E2').Value = '  +0.03%  '

$ws.Range('D3').Value = '1.571.70'
$ws.Range('E3').Value = '  +0.24%  '

$ws.Range('E4').Value = '  +0.06%  '

$ws.Range('D5').Value = '''1.001'

$ws.Range('D6').Value = '''291.59'
$ws.Range('E6').Value = '  +0.34%  '

$ws.Range('D7').Value = '''0.3759'
$ws.Range('E7').Value = '  +2.08%  '

$ws.Range('D8').Value = '''50.00'
$ws.Range('E8').Value = '  +1.09%  '

$ws.Range('D9').Value = '''0.3422'
$ws.Range('E9').Value = '  +0.87%  '

$ws.Range('D10').Value = '''0.07650'
$ws.Range('E10').Value = '  +0.50%  '

$ws.Range('D11').Value = '''1.152'
$ws.Range('E11').Value = '  -1.61%  '

$ws.Range('E12').Value = '  +0.07%  '

$ws.Range('D13').Value = '''21.20'
$ws.Range('E13').Value = '  -0.10%  '

$ws.Range('D14').Value = '''6.008'
$ws.Range('E14').Value = '  -0.78%  '

$ws.Range('D15').Value = '''6.941'
$ws.Range('E15').Value = '  +0.62%  '

$ws.Range('D16').Value = '1.570.32'
$ws.Range('E16').Value = '  -0.56%  '

$ws.Range('D17').Value = '''0.00001131'
$ws.Range('E17').Value = '  -0.44%  '

$ws.Range('D18').Value = '''89.92'
$ws.Range('E18').Value = '  +0.91%  '

$ws.Range('D19').Value = '''0.06760'
$ws.Range('E19').Value = '  -0.01%  '

$ws.Range('E20').Value = '  -0.06%  '

$ws.Range('D21').Value = '''16.79'
$ws.Range('E21').Value = '  +1.63%  '

$ws.Range('D22').Value = '''6.208'
$ws.Range('E22').Value = '  -0.42%  '

$ws.Range('D23').Value = '''12.00'
$ws.Range('E23').Value = '  -0.41%  '

$ws.Range('D24').Value = '22.383.60'
$ws.Range('E24').Value = '  -0.10%  '

$ws.Range('D25').Value = '''2.397'
$ws.Range('E25').Value = '  +0.42%  '

$ws.Range('D26').Value = '''2.680'
$ws.Range('E26').Value = '  -10.23%  '

$ws.Range('D27').Value = '''20.23'
$ws.Range('E27').Value = '  +1.64%  '

$ws.Range('D28').Value = '''147.19'
$ws.Range('E28').Value = '  +1.04%  '

$ws.Range('D29').Value = '''5.028'
$ws.Range('E29').Value = '  +1.55%  '

$ws.Range('D30').Value = '''126.25'
$ws.Range('E30').Value = '  +0.75%  '

$ws.Range('D31').Value = '1.743.70'
$ws.Range('E31').Value = '  -0.66%  '

$ws.Range('D32').Value = '''6.155'
$ws.Range('E32').Value = '  -1.68%  '

$ws.Range('E33').Value = '  +1.00%  '

$ws.Range('D34').Value = '''0.9814'
$ws.Range('E34').Value = '  -5.53%  '

$ws.Range('D35').Value = '''9.912'
$ws.Range('E35').Value = '  -3.67%  '

$ws.Range('D36').Value = '''0.08507'
$ws.Range('E36').Value = '  +0.53%  '

$ws.Range('D37').Value = '''0.02550'
$ws.Range('E37').Value = '  +0.19%  '

$ws.Range('D38').Value = '''1.383'
$ws.Range('E38').Value = '  +11.06%  '

$ws.Range('D39').Value = '''0.2318'
$ws.Range('E39').Value = '  -0.46%  '

$ws.Range('D40').Value = '''0.06561'
$ws.Range('E40').Value = '  +0.11%  '

$ws.Range('D41').Value = '''5.428'
$ws.Range('E41').Value = '  -1.98%  '

$ws.Range('D42').Value = '''0.6398'
$ws.Range('E42').Value = '  +0.29%  '

$ws.Range('D43').Value = '''11.46'
$ws.Range('E43').Value = '  -3.22%  '

$ws.Range('E44').Value = '  +0.05%  '

$ws.Range('D45').Value = '''14.07'
$ws.Range('E45').Value = '  -2.17%  '

$ws.Range('D46').Value = '''3.780'

$ws.Range('D47').Value = '''0.5975'
$ws.Range('E47').Value = '  -0.65%  '

$ws.Range('D48').Value = '''1.293'
$ws.Range('E48').Value = '  +2.28%  '

$ws.Range('D49').Value = '''2.090'
$ws.Range('E49').Value = '  -2.08%  '

$ws.Range('D50').Value = '''125.47'
$ws.Range('E50').Value = '  +1.72%  '

$ws.Range('E51').Value = '  +0.61%  '
